# Update timestamps on the "data" sheet (column F, rows 2-22)
$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

$timestamps = @(
    "2021-10-05 14:22:20.825012",
    "2021-10-05 14:22:20.825017",
    "2021-10-05 14:22:20.825020",
    "2021-10-05 14:22:20.825022",
    "2021-10-05 14:22:20.825024",
    "2021-10-05 14:22:20.825026",
    "2021-10-05 14:22:20.825028",
    "2021-10-05 14:22:20.825030",
    "2021-10-05 14:22:20.825032",
    "2021-10-05 14:22:20.825034",
    "2021-10-05 14:22:20.825036",
    "2021-10-05 14:22:20.825038",
    "2021-10-05 14:22:20.825040",
    "2021-10-05 14:22:20.825042",
    "2021-10-05 14:22:20.825044",
    "2021-10-05 14:22:20.825045",
    "2021-10-05 14:22:20.825048",
    "2021-10-05 14:22:20.825050",
    "2021-10-05 14:22:20.825051",
    "2021-10-05 14:22:20.825053",
    "2021-10-05 14:22:20.825055"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $timestamps[$i]
}

# Add the "metadata" sheet as a new second tab, directly after "data". We
# build it by copying "data" (so it inherits the same sheetPr/outline props
# and page margins) and then clearing it down to the rows/columns we need -
# this keeps the original header style (s="1") intact with no new styles.
$data.Copy($null, $data)
$meta = $wb.Worksheets.Item(2)
$meta.Name = "metadata"

# Drop the copied data rows below the header/first row; only 2 rows remain.
$meta.Rows("3:22").Delete()

$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

# G1 is a brand new header column - copy the header format from F1 so it
# reuses the same style index instead of registering a new one.
$meta.Range("F1").Copy()
$meta.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 2).Value = "Progressive cardiac conduction disease"
$meta.Cells.Item(2, 3).Value = 506
$meta.Cells.Item(2, 5).Value = "2021-03-02T15:56:22.257376Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:22:20.822590"
$meta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/506/?format=json"

# data_version ("1.37") must be stored as TEXT, not a number, to match the
# source data. A plain .Value assignment of a numeric-looking string gets
# auto-coerced to a number (same as typing it into Excel), so instead we
# compute it with TEXT() and then collapse the formula down to its cached
# string result via a self value-paste (keeps the default, unstyled cell).
$verCell = $meta.Cells.Item(2, 4)
$verCell.Formula = '=TEXT(1.37,"0.00")'
$verCell.Copy()
$verCell.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# Keep "data" as the active sheet/tab
$data.Activate() | Out-Null
$data.Range("A1").Select() | Out-Null
